# Weekly fruit/vegetable price update: a new daily price record is
# inserted as row 37 ("Región Metropolitana", fecha 44546), pushing all
# the following rows (previously 37-83) down by one (now 38-84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37; Excel shifts rows 37..83 down to 38..84 and
# carries the existing per-column formatting (e.g. the date style on D).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new price record.
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44546
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112032
$ws.Range("G37").Value = "Zapallo italiano"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 250
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 9000
$ws.Range("M37").Value = 8400
$ws.Range("N37").Value = "$/caja 60 unidades"
$ws.Range("O37").Value = "Región Metropolitana"
$ws.Range("P37").Value = 140
$ws.Range("Q37").Value = 60
$ws.Range("R37").Value = "Hortaliza"
